$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Groups of rows whose B:AC content gets cyclically rotated.
# For a group [r0, r1, ..., rn], the new content of r0 becomes the old
# content of r1, new r1 becomes old r2, ..., new rn becomes old r0.
$groups = @(
    @(327, 328),
    @(354, 355),
    @(362, 363, 364),
    @(372, 373),
    @(395, 396),
    @(405, 406),
    @(417, 418),
    @(422, 423),
    @(424, 426),
    @(425, 427),
    @(453, 454),
    @(455, 457),
    @(506, 507),
    @(510, 511),
    @(588, 589),
    @(596, 597),
    @(680, 682, 681)
)

foreach ($group in $groups) {
    # Snapshot the "before" values of B:AC for every row in this group.
    $snapshots = @()
    foreach ($r in $group) {
        $snapshots += ,($ws.Range("B" + $r + ":AC" + $r).Value())
    }

    $count = $group.Length
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $group[$i]
        $srcIndex = ($i + 1) % $count
        $ws.Range("B" + $destRow + ":AC" + $destRow).Value = $snapshots[$srcIndex]
    }
}
